$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1): N1, O1, P1 ---
$ws.Range("N1").Value = "ScmEmails"
$ws.Range("O1").Value = "CustomsReleaseEmails"
$ws.Range("P1").Value = "BillingEmails"
$ws.Range("N1:P1").Font.Bold = $true

# --- New data cells (row 2): N2, O2, P2 ---
$ws.Range("N2").Value = "topite@mailinator.com; julio.wei@1000shores.com; yiwudoc@1000shores.com"
$ws.Range("P2").Value = "yhaglcwuliu@163.com; 295362722@qq.com; 2590375680@qq.com; 13101609583@163.com"

# O2 gets a hyperlink (and picks up the built-in Hyperlink style automatically)
$ws.Hyperlinks.Add($ws.Range("O2"), "mailto:chenjiao@lionifreight.com;kliya@lionifreight.com;yezhixuan@lionifreight.com", [Type]::Missing, [Type]::Missing, "chenjiao@lionifreight.com;kliya@lionifreight.com;yezhixuan@lionifreight.com")

# --- Add a tooltip (ScreenTip) to the existing D2 hyperlink ---
foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Address() -eq '$D$2') {
        $h.ScreenTip = "mailto:topite@mailinator.com"
    }
}

# --- Update selection / view ---
$excel.Goto($ws.Range("K6"), $true)
